$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 199 entirely (the stray theoretical-line value 4942.4592),
# shifting every row below it up by one, per the corrected "velocidad" list.
$ws.Rows.Item(199).Delete()

# Update the saved view state: scroll position and active cell selection.
$excel.ActiveWindow.ScrollRow = 193
$ws.Range("F210").Select()
